# Applies the "Nombreuses informations supplémentaires récupérées." commit:
#  - "exemplars" property renamed to "copies" (Bibliographic Item sheet, rows 17-21)
#  - "directUrl" property renamed to "url" (Bibliographic Item sheet, row 22)
#  - a couple of newly-recovered "Disponible" flags added (F22, F24)
#  - active sheet/selection moved from "Result Set" back to "Bibliographic Item"

$wb = $excel.ActiveWorkbook
$wsBib = $wb.Worksheets.Item(1)   # "Bibliographic Item"

# Rename "exemplars" -> "copies" for the five rows describing that property.
$wsBib.Range("A17").Value = "copies"
$wsBib.Range("A18").Value = "copies"
$wsBib.Range("A19").Value = "copies"
$wsBib.Range("A20").Value = "copies"
$wsBib.Range("A21").Value = "copies"

# Rename "directUrl" -> "url" for the directAccesses property row.
$wsBib.Range("B22").Value = "url"

# Newly recovered data points now mark these columns as available too.
$wsBib.Range("F22").Value = "Disponible"
$wsBib.Range("F24").Value = "Disponible"

# Switch the active tab/selection back to the "Bibliographic Item" sheet.
$wsBib.Activate()
$wsBib.Range("A24").Select()
